$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell address, new text value (applied as literal text via
# a leading apostrophe so Excel does not reinterpret numeric-looking / percent-
# looking strings as numbers; the quote-prefix style bump is then undone by
# resetting the cell to the 'Normal' style, restoring the original (unstyled) xf.
$updates = @(
    @{Cell='D2'; Value='275.32'}
    @{Cell='E2'; Value='-1.23%'}
    @{Cell='G2'; Value='12'}
    @{Cell='D3'; Value='26.58'}
    @{Cell='E3'; Value='-3.07%'}
    @{Cell='G3'; Value='12'}
    @{Cell='D4'; Value='4.867'}
    @{Cell='E4'; Value='1.59%'}
    @{Cell='G4'; Value='12'}
    @{Cell='D5'; Value='0.06311'}
    @{Cell='E5'; Value='-0.45%'}
    @{Cell='G5'; Value='12'}
    @{Cell='D6'; Value='6.931'}
    @{Cell='E6'; Value='-0.10%'}
    @{Cell='G6'; Value='12'}
    @{Cell='D7'; Value='1.293'}
    @{Cell='E7'; Value='35.11%'}
    @{Cell='G7'; Value='12'}
    @{Cell='E8'; Value='-1.15%'}
    @{Cell='G8'; Value='12'}
    @{Cell='D9'; Value='0.1533'}
    @{Cell='E9'; Value='4.28%'}
    @{Cell='G9'; Value='12'}
    @{Cell='D10'; Value='0.05063'}
    @{Cell='E10'; Value='-1.30%'}
    @{Cell='G10'; Value='12'}
    @{Cell='D11'; Value='0.07407'}
    @{Cell='E11'; Value='1.69%'}
    @{Cell='G11'; Value='12'}
    @{Cell='D12'; Value='0.02898'}
    @{Cell='E12'; Value='-7.64%'}
    @{Cell='G12'; Value='12'}
    @{Cell='D13'; Value='0.09064'}
    @{Cell='E13'; Value='-0.01%'}
    @{Cell='G13'; Value='12'}
    @{Cell='D14'; Value='0.001600'}
    @{Cell='E14'; Value='2.74%'}
    @{Cell='G14'; Value='12'}
    @{Cell='D15'; Value='0.0006352'}
    @{Cell='E15'; Value='0.93%'}
    @{Cell='G15'; Value='12'}
    @{Cell='D16'; Value='0.006076'}
    @{Cell='E16'; Value='4.60%'}
    @{Cell='G16'; Value='12'}
    @{Cell='D17'; Value='3.455'}
    @{Cell='E17'; Value='0.10%'}
    @{Cell='G17'; Value='12'}
    @{Cell='D18'; Value='3.310'}
    @{Cell='E18'; Value='-2.28%'}
    @{Cell='G18'; Value='12'}
    @{Cell='G19'; Value='12'}
    @{Cell='E20'; Value='-0.20%'}
    @{Cell='G20'; Value='12'}
    @{Cell='D21'; Value='0.1313'}
    @{Cell='E21'; Value='-1.40%'}
    @{Cell='G21'; Value='12'}
    @{Cell='D22'; Value='3.905'}
    @{Cell='E22'; Value='1.23%'}
    @{Cell='G22'; Value='12'}
    @{Cell='D23'; Value='0.04397'}
    @{Cell='E23'; Value='1.97%'}
    @{Cell='G23'; Value='12'}
    @{Cell='E24'; Value='-0.56%'}
    @{Cell='G24'; Value='12'}
    @{Cell='E25'; Value='-1.78%'}
    @{Cell='G25'; Value='12'}
    @{Cell='D26'; Value='0.0001201'}
    @{Cell='G26'; Value='12'}
    @{Cell='D27'; Value='0.0001659'}
    @{Cell='E27'; Value='-1.84%'}
    @{Cell='G27'; Value='12'}
    @{Cell='G28'; Value='12'}
    @{Cell='G29'; Value='12'}
    @{Cell='G30'; Value='12'}
    @{Cell='G31'; Value='12'}
    @{Cell='G32'; Value='12'}
    @{Cell='G33'; Value='12'}
    @{Cell='G34'; Value='12'}
    @{Cell='G35'; Value='12'}
    @{Cell='G36'; Value='12'}
    @{Cell='G37'; Value='12'}
    @{Cell='G38'; Value='12'}
    @{Cell='G39'; Value='12'}
    @{Cell='D40'; Value='0.04077'}
    @{Cell='E40'; Value='0.02%'}
    @{Cell='G40'; Value='12'}
    @{Cell='D41'; Value='0.007004'}
    @{Cell='E41'; Value='4.45%'}
    @{Cell='G41'; Value='12'}
    @{Cell='E42'; Value='0.85%'}
    @{Cell='G42'; Value='12'}
    @{Cell='D43'; Value='0.002142'}
    @{Cell='E43'; Value='-2.65%'}
    @{Cell='G43'; Value='12'}
    @{Cell='D44'; Value='0.01119'}
    @{Cell='E44'; Value='-10.40%'}
    @{Cell='G44'; Value='12'}
    @{Cell='D45'; Value='0.00005208'}
    @{Cell='E45'; Value='-0.03%'}
    @{Cell='G45'; Value='12'}
    @{Cell='B46'; Value='CoinbaseStockToken'}
    @{Cell='C46'; Value='https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'}
    @{Cell='D46'; Value='0.02000'}
    @{Cell='E46'; Value='-11.20%'}
    @{Cell='G46'; Value='12'}
    @{Cell='B47'; Value='BOLO'}
    @{Cell='C47'; Value='https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'}
    @{Cell='D47'; Value='1.490'}
    @{Cell='E47'; Value='-37.38%'}
    @{Cell='G47'; Value='12'}
    @{Cell='G48'; Value='12'}
    @{Cell='G49'; Value='12'}
    @{Cell='G50'; Value='12'}
    @{Cell='G51'; Value='12'}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.Value = "'" + $u.Value
    $rng.Style = 'Normal'
}
